$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-03 Friday" "2025-10-04 Saturday"

Replace-Text "38÷9=" "84÷2="
Replace-Text "23÷8=" "12÷5="
Replace-Text "44÷4=" "84÷5="
Replace-Text "65÷2=" "80÷2="
Replace-Text "18÷5=" "19÷7="

Replace-Text "88÷7=" "43÷8="
Replace-Text "31÷4=" "12÷7="
Replace-Text "46÷4=" "91÷8="
Replace-Text "30÷5=" "77÷4="
Replace-Text "29÷9=" "81÷3="

Replace-Text "73÷9=" "28÷7="
Replace-Text "77÷3=" "38÷5="
Replace-Text "73÷7=" "46÷6="
Replace-Text "87÷8=" "83÷2="
Replace-Text "82÷7=" "80÷9="

Replace-Text "78÷6=" "81÷8="
Replace-Text "20÷8=" "51÷7="
Replace-Text "74÷6=" "26÷2="
Replace-Text "90÷2=" "73÷5="
Replace-Text "76÷2=" "38÷7="

Replace-Text "63÷4=" "99÷8="
Replace-Text "29÷4=" "87÷4="
Replace-Text "97÷4=" "72÷3="
Replace-Text "40÷6=" "35÷7="
Replace-Text "81÷5=" "99÷2="
